# Update gh-pages to output generated at 456a3b4
# Applies the "F column" (想去人数 / wants-to-go count) refresh plus the
# newly scraped "九江·动漫畅想" event row to the 展览 (sheet1) and
# 全部类型 (sheet4) tabs, and a single F-column refresh on 演出 (sheet2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------
# 1) 展览 (sheet1): refresh "want to go" counts (column F) for existing rows
# ---------------------------------------------------------------------
$ws1.Cells.Item(4,6).Value = 4809
$ws1.Cells.Item(6,6).Value = 156
$ws1.Cells.Item(8,6).Value = 108
$ws1.Cells.Item(9,6).Value = 0
$ws1.Cells.Item(10,6).Value = 753
$ws1.Cells.Item(11,6).Value = 0
$ws1.Cells.Item(15,6).Value = 185
$ws1.Cells.Item(16,6).Value = 0
$ws1.Cells.Item(18,6).Value = 112
$ws1.Cells.Item(19,6).Value = 3996
$ws1.Cells.Item(20,6).Value = 0
$ws1.Cells.Item(23,6).Value = 0
$ws1.Cells.Item(24,6).Value = 0
$ws1.Cells.Item(26,6).Value = 0
$ws1.Cells.Item(27,6).Value = 0
$ws1.Cells.Item(29,6).Value = 0
$ws1.Cells.Item(31,6).Value = 0
$ws1.Cells.Item(33,6).Value = 0
$ws1.Cells.Item(34,6).Value = 287
$ws1.Cells.Item(36,6).Value = 0
$ws1.Cells.Item(37,6).Value = 175
$ws1.Cells.Item(38,6).Value = 0
$ws1.Cells.Item(39,6).Value = 0
$ws1.Cells.Item(40,6).Value = 969
$ws1.Cells.Item(42,6).Value = 70

# Insert the new "九江·动漫畅想" event as row 46 (pushes the old rows
# 46 "上饶·..." -> 47 and 47 "南昌·..." -> 48).
$ws1.Rows.Item(46).Insert()

# Match the bold/bordered "index" column style used throughout column A.
$ws1.Cells.Item(45,1).Copy()
$ws1.Cells.Item(46,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Cells.Item(46,1).Value = 45
$ws1.Cells.Item(46,2).NumberFormat = "@"
$ws1.Cells.Item(46,2).Value = "2024-08-11"
$ws1.Cells.Item(46,3).Value = "九江·动漫畅想"
$ws1.Cells.Item(46,4).Value = "十里大道202号（十里大道与地质路交汇处） 山水国际大酒店"
$ws1.Cells.Item(46,5).Value = "2024.08.11 08:00-08.11 20:00"
$ws1.Cells.Item(46,6).Value = 0
$ws1.Cells.Item(46,7).Value = 50
$ws1.Cells.Item(46,8).Value = "https://show.bilibili.com/platform/detail.html?id=89029"
$ws1.Cells.Item(46,9).Value = "//i1.hdslb.com/bfs/openplatform/202407/NLseFxU51720516345581.png"

# Renumber the pushed-down rows' index column (A is a plain 0-based counter).
$ws1.Cells.Item(47,1).Value = 46
$ws1.Cells.Item(48,1).Value = 47

# ---------------------------------------------------------------------
# 2) 演出 (sheet2): refresh "want to go" count for row 2
# ---------------------------------------------------------------------
$ws2.Cells.Item(2,6).Value = 0

# ---------------------------------------------------------------------
# 3) 全部类型 (sheet4): refresh "want to go" counts (column F)
# ---------------------------------------------------------------------
$ws4.Cells.Item(2,6).Value = 32
$ws4.Cells.Item(4,6).Value = 4809
$ws4.Cells.Item(6,6).Value = 0
$ws4.Cells.Item(9,6).Value = 108
$ws4.Cells.Item(10,6).Value = 94
$ws4.Cells.Item(11,6).Value = 0
$ws4.Cells.Item(13,6).Value = 1180
$ws4.Cells.Item(14,6).Value = 0
$ws4.Cells.Item(16,6).Value = 0
$ws4.Cells.Item(17,6).Value = 0
$ws4.Cells.Item(18,6).Value = 151
$ws4.Cells.Item(20,6).Value = 3996
$ws4.Cells.Item(21,6).Value = 6326
$ws4.Cells.Item(22,6).Value = 0
$ws4.Cells.Item(23,6).Value = 0
$ws4.Cells.Item(25,6).Value = 540
$ws4.Cells.Item(26,6).Value = 48
$ws4.Cells.Item(27,6).Value = 3972
$ws4.Cells.Item(31,6).Value = 0
$ws4.Cells.Item(33,6).Value = 533
$ws4.Cells.Item(35,6).Value = 0
$ws4.Cells.Item(38,6).Value = 175
$ws4.Cells.Item(40,6).Value = 0
$ws4.Cells.Item(41,6).Value = 969
$ws4.Cells.Item(43,6).Value = 0
$ws4.Cells.Item(45,6).Value = 495
$ws4.Cells.Item(47,6).Value = 0

# Insert the new "九江·动漫畅想" event as row 47 (pushes the old rows
# 47 "上饶·..." -> 48 and 48 "南昌·..." -> 49).
$ws4.Rows.Item(47).Insert()

$ws4.Cells.Item(46,1).Copy()
$ws4.Cells.Item(47,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws4.Cells.Item(47,1).Value = 46
$ws4.Cells.Item(47,2).NumberFormat = "@"
$ws4.Cells.Item(47,2).Value = "2024-08-11"
$ws4.Cells.Item(47,3).Value = "九江·动漫畅想"
$ws4.Cells.Item(47,4).Value = "十里大道202号（十里大道与地质路交汇处） 山水国际大酒店"
$ws4.Cells.Item(47,5).Value = "2024.08.11 08:00-08.11 20:00"
$ws4.Cells.Item(47,6).Value = 0
$ws4.Cells.Item(47,7).Value = 50
$ws4.Cells.Item(47,8).Value = "https://show.bilibili.com/platform/detail.html?id=89029"
$ws4.Cells.Item(47,9).Value = "//i1.hdslb.com/bfs/openplatform/202407/NLseFxU51720516345581.png"

# Renumber the pushed-down rows' index column.
$ws4.Cells.Item(48,1).Value = 47
$ws4.Cells.Item(49,1).Value = 48

Write-Output "Edit applied"
